# Auto-generated edit script: apply Kraken_Profits value updates per sheet
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @{
    "H13" = 1344
    "I13" = 1310
    "J13" = 1378
    "K13" = 1310
    "L13" = 1378
    "M13" = -1141
    "N13" = -1716
    "H98" = 1487.091
    "I98" = 1487.091
    "K98" = 1487.091
    "M98" = 10.90900000000011
    "H122" = 1487.091
    "I122" = 1487.091
    "K122" = 4461.272999999999
    "M122" = -2011.272999999999
    "H137" = 4460.357
    "J137" = 5589.6
    "L137" = 16768.8
    "N137" = -21868.8
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @{
    "H11" = 5388.2
    "I11" = 4471
    "J11" = 5999.6665
    "K11" = 4471
    "L11" = 5999.6665
    "M11" = -4327
    "N11" = -6287.6665
    "H12" = 898
    "J12" = 975
    "L12" = 975
    "N12" = -1321
    "H15" = 3000
    "J15" = 3000
    "L15" = 3000
    "N15" = -3700
    "H17" = 500
    "I17" = 500
    "J17" = 0
    "K17" = 500
    "L17" = 0
    "M17" = -327
    "H21" = 0
    "I21" = 0
    "K21" = 0
    "H22" = 0
    "I22" = 0
    "J22" = 0
    "K22" = 0
    "L22" = 0
    "H23" = 12500
    "I23" = 10000
    "J23" = 15000
    "K23" = 10000
    "L23" = 15000
    "M23" = -9741
    "N23" = -15518
    "H32" = 2384.05
    "I32" = 2599.0588
    "K32" = 2599.0588
    "M32" = -2312.0588
    "H88" = 1560
    "J88" = 1560
    "L88" = 1560
    "N88" = -2372
    "H91" = 1560
    "J91" = 1560
    "L91" = 1560
    "N91" = -4368
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$clears = @("N17", "M21", "M22", "N22")
foreach ($addr in $clears) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @{
    "H5" = 1000
    "J5" = 1000
    "L5" = 1000
    "N5" = -1226
    "H7" = 6666967
    "I7" = 10000200
    "K7" = 10000200
    "M7" = -10000087
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @{
    "H5" = 1106.3334
    "J5" = 1627
    "L5" = 1627
    "N5" = -1851
    "H11" = 853
    "I11" = 1500
    "K11" = 1500
    "M11" = -1360
    "H13" = 2500354.5
    "I13" = 3333639.2
    "J13" = 500
    "K13" = 3333639.2
    "L13" = 500
    "M13" = -3333500.2
    "N13" = -778
    "H41" = 20000
    "I41" = 0
    "K41" = 0
    "H58" = 0
    "I58" = 0
    "J58" = 0
    "K58" = 0
    "L58" = 0
    "H86" = 4516.6665
    "J86" = 3850
    "L86" = 3850
    "N86" = -6096
    "H89" = 4516.6665
    "J89" = 3850
    "L89" = 19250
    "N89" = -30482
    "H136" = 0
    "I136" = 0
    "J136" = 0
    "K136" = 0
    "L136" = 0
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$clears = @("M41", "M58", "N58", "M136", "N136")
foreach ($addr in $clears) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @{
    "H23" = 1614.2858
    "I23" = 2000
    "J23" = 1460
    "K23" = 6000
    "L23" = 4380
    "M23" = -5765
    "N23" = -4850
    "H80" = 5800
    "H83" = 5800
    "H103" = 1796.5714
    "J103" = 2206.5
    "L103" = 6619.5
    "N103" = -8377.5
    "H117" = 1031.6666
    "J117" = 1031.6666
    "L117" = 3094.9998
    "N117" = -9978.9998
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @{
    "H136" = 38794.555
    "J136" = 38794.555
    "L136" = 116383.665
    "N136" = -121483.665
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @{
    "H22" = 2211.7188
    "I22" = 1920.6522
    "J22" = 2955.5557
    "K22" = 1920.6522
    "L22" = 2955.5557
    "M22" = -1625.6522
    "N22" = -3545.5557
    "H27" = 2211.7188
    "I27" = 1920.6522
    "J27" = 2955.5557
    "K27" = 1920.6522
    "L27" = 2955.5557
    "M27" = -1813.6522
    "N27" = -3169.5557
    "H46" = 3778.4614
    "I46" = 0
    "J46" = 3778.4614
    "K46" = 0
    "L46" = 3778.4614
    "N46" = -4154.4614
    "H93" = 4886.2856
    "I93" = 4867.5
    "K93" = 4867.5
    "M93" = -3619.5
    "H100" = 9988
    "I100" = 3178.4
    "K100" = 3178.4
    "M100" = -2637.4
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
$clears = @("M46")
foreach ($addr in $clears) {
    $ws.Range($addr).ClearContents()
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @{
    "H3" = 12650
    "I3" = 18002
    "K3" = 18002
    "M3" = -17888
    "H8" = 1000
    "J8" = 1000
    "L8" = 1000
    "N8" = -1280
    "H11" = 6625.75
    "I11" = 3004
    "K11" = 3004
    "M11" = -2862
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
